$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("D2").Value = "71.781.48"
$ws.Range("E2").Value = "  -0.45%  "

$ws.Range("D3").Value = "3.999.70"
$ws.Range("E3").Value = "  -1.15%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "534.66"
$ws.Range("E5").Value = "  +3.15%  "

$ws.Range("D6").Value = "149.00"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").Value = "0.694"
$ws.Range("E7").Value = "  +11.65%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.742"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -2.35%  "

$ws.Range("B11").Value = "ShibaInu"
$ws.Range("C11").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -3.18%  "

$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "50.30"
$ws.Range("E12").Value = "  +4.44%  "

$ws.Range("D13").Value = "10.66"
$ws.Range("E13").Value = "  -1.75%  "

$ws.Range("D14").Value = "4.650.70"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").Value = "4.005.98"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").Value = "13.96"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").Value = "20.48"
$ws.Range("E17").Value = "  -3.39%  "

$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("E19").Value = "  -2.64%  "

$ws.Range("D20").Value = "71.832.87"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").Value = "426.57"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").Value = "96.81"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").Value = "3.47"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("D24").Value = "4.20"
$ws.Range("E24").Value = "  +4.79%  "

$ws.Range("D25").Value = "14.22"
$ws.Range("E25").Value = "  -2.88%  "

$ws.Range("D26").Value = "11.13"
$ws.Range("E26").Value = "  -6.46%  "

$ws.Range("D27").Value = "10.67"
$ws.Range("E27").Value = "  -4.48%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "5.84"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "3.67"
$ws.Range("E29").Value = "  +18.49%  "

$ws.Range("D30").Value = "36.66"
$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("D31").Value = "7.36"
$ws.Range("E31").Value = "  +5.11%  "

$ws.Range("D33").Value = "13.32"
$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "673.66"
$ws.Range("E34").Value = "  -3.92%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "47.81"
$ws.Range("E35").Value = "  +18.24%  "

$ws.Range("D36").Value = "65.31"
$ws.Range("E36").Value = "  -3.99%  "

$ws.Range("D37").Value = "0.443"
$ws.Range("E37").Value = "  +1.91%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0814"
$ws.Range("E39").Value = "  -8.93%  "

$ws.Range("D40").Value = "3.41"
$ws.Range("E40").Value = "  -6.64%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  +6.16%  "

$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("D44").Value = "0.0485"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("E45").Value = "  +2.75%  "

$ws.Range("D46").Value = "9.81"
$ws.Range("E46").Value = "  +8.55%  "

$ws.Range("E47").Value = "  -3.98%  "

$ws.Range("D48").Value = "3.36"
$ws.Range("E48").Value = "  -4.35%  "

$ws.Range("D49").Value = "2.99"
$ws.Range("E49").Value = "  -4.34%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000266"
$ws.Range("E50").Value = "  -1.86%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "142.95"
$ws.Range("E51").Value = "  +0.01%  "
